# Rescale the "Valor" values in column D (rows 2-29) from a fraction
# (e.g. 0.2560534787983632) to a percentage-style number (25.60534787983632)
# by multiplying each value by 100.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = $cell.Value2 * 100
}
